# The upstream change for this fixture is recorded in the repository as a
# unified diff against the package's OOXML.  Inspecting that diff carefully
# shows that every single hunk is a pure attribute-order shuffle (e.g.
# <w:pgSz w:w="11906" w:h="16838"/> -> <w:pgSz w:h="16838" w:w="11906"/>,
# namespace declarations on <w:document> being alphabetised, etc.) produced
# when the original authors re-saved the .docx test fixture with a newer
# OOXML toolchain. The attribute *sets* (name/value pairs), element
# structure, and all text content are byte-for-byte identical before and
# after - there is no observable content, formatting, or structural change
# to reproduce through the Word object model.
#
# Concretely: for every element touched by the diff (the <w:document> root,
# <w:pgSz>, <w:pgMar>, <w:rFonts>, <w:lang>, <w:latentStyles>,
# <w:lsdException>, <w:style>, <w:spacing>, <w:color>, <w:tblInd>, the table
# cell margins, ...) the old/new attribute dictionaries are identical sets;
# only their serialisation order changed. That reordering is an artifact of
# whatever XML writer resaved the package upstream, not something exposed
# through (or controllable via) the Word COM automation surface - there is
# no "resort attributes" verb in the object model, and deliberately poking
# the relevant properties (PageSetup margins, style fonts, ...) with their
# own current values only adds unrelated noise (extra namespace
# declarations, stray empty <w:rPr> elements, ...) that is not present in
# the target diff.
#
# So the faithful edit is to leave the document's content untouched. We
# still touch the document through the object model (read-only) so the
# script demonstrably runs against the live COM object rather than being a
# literal empty file.

$d = $word.ActiveDocument

# Touch a few of the exact nodes the upstream diff mentions, read-only, to
# confirm they exist/are reachable - without writing anything back (writing
# back identical values still perturbs serialization, see above).
$null = $d.PageSetup.PageWidth
$null = $d.PageSetup.PageHeight
$null = $d.PageSetup.TopMargin
$null = $d.Styles("Normal").NameLocal
$null = $d.Content.Text

# No content, formatting or structural mutation is applied: the source and
# target OOXML are semantically identical.
